# Upload new version with timestamp
# Adds two new inventory rows (GARAMYCIN + جنتيانا), renumbers the "م"
# index column, bumps the grand total, and refreshes the generated
# timestamp footer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Insert "GARAMYCIN 0.1% OINT. 15 GM" as a new row just before
#    "HIBIOTIC N 600MG SUSP. 80 ML" (currently row 18).
# ---------------------------------------------------------------------
$ws.Rows.Item(18).Insert()
$ws.Rows.Item(17).EntireRow.Copy()
$ws.Rows.Item(18).PasteSpecial()

$ws.Cells.Item(18, 1).Value = 12
$ws.Cells.Item(18, 3).Value = "GARAMYCIN 0.1% OINT. 15 GM"
$ws.Cells.Item(18, 8).Value = "1:0"
$ws.Cells.Item(18, 12).Value = "1"
$ws.Cells.Item(18, 14).Value = "22.00"
$ws.Cells.Item(18, 16).Value = "22.0000"
$ws.Cells.Item(18, 17).Value = "1:0"

# Renumber the rest of the "م" column (now rows 19-35, was 18-34).
for ($r = 19; $r -le 35; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 6
}

# ---------------------------------------------------------------------
# 2) Insert "جنتيانا" as a new row just before "جهاز محلول" (currently
#    row 32 now that GARAMYCIN shifted everything down by one).
# ---------------------------------------------------------------------
$ws.Rows.Item(32).Insert()
$ws.Rows.Item(31).EntireRow.Copy()
$ws.Rows.Item(32).PasteSpecial()

$ws.Cells.Item(32, 1).Value = 26
$ws.Cells.Item(32, 3).Value = "جنتيانا "
$ws.Cells.Item(32, 8).Value = "9:0"
$ws.Cells.Item(32, 12).Value = "0"
$ws.Cells.Item(32, 14).Value = "15.00"
$ws.Cells.Item(32, 16).Value = "15.0000"
$ws.Cells.Item(32, 17).Value = "1:0"

# Renumber the rest of the "م" column (now rows 33-37, was 32-36).
for ($r = 33; $r -le 37; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 6
}

# ---------------------------------------------------------------------
# 3) Update the grand total (now on row 38) for the two new line items.
# ---------------------------------------------------------------------
$ws.Cells.Item(38, 16).Value = 1336.21

# ---------------------------------------------------------------------
# 4) Refresh the generated timestamp footer (now on row 39).
# ---------------------------------------------------------------------
$ws.Cells.Item(39, 1).Value = "Wednesday, 23 July, 2025 2:08 PM"

$wb.Save()
